$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B28").Value = "https://drive.google.com/file/d/12g4aTrsdVaLsj6B9soJ1t-MSPUay6PJd/view?usp=drive_link"
$ws.Hyperlinks.Add($ws.Range("B28"), "https://drive.google.com/file/d/12g4aTrsdVaLsj6B9soJ1t-MSPUay6PJd/view?usp=drive_link")
$ws.Range("B28").Style = "Hyperlink"
Write-Output "done"
